$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. "external" list (column J): add terminate(programName) after tail(id,file) ---
$ws.Cells.Item(6, 10).Value = "terminate(programName)"

# --- 2. "base" list (column F): remove clearVariables(variables) (row 19),
#        shifting rows 20-40 up by one (F19:F39), clearing the now-empty F40 ---
$baseValues = @(
    "failImmediate(text)",
    "incrementChar(var,amount,config)",
    "macro(file,sheet,name)",
    "outputToCloud(resource)",
    "prependText(var,prependWith)",
    "repeatUntil(steps,maxWaitMs)",
    "save(var,value)",
    "saveCount(text,regex,saveVar)",
    "saveMatches(text,regex,saveVar)",
    "saveReplace(text,regex,replace,saveVar)",
    "saveVariablesByPrefix(var,prefix)",
    "saveVariablesByRegex(var,regex)",
    "section(steps)",
    "split(text,delim,saveVar)",
    "startRecording()",
    "stopRecording()",
    "substringAfter(text,delim,saveVar)",
    "substringBefore(text,delim,saveVar)",
    "substringBetween(text,start,end,saveVar)",
    "verbose(text)",
    "waitFor(waitMs)"
)
for ($i = 0; $i -lt $baseValues.Length; $i++) {
    $ws.Cells.Item(19 + $i, 6).Value = $baseValues[$i]
}
$ws.Cells.Item(40, 6).ClearContents()

# --- 3. "web" list (column Z): insert two new commands before saveTableAsCsv(...)
#        (old rows 99-135 move down to 101-137) ---
$webValues = @(
    "saveTableAsCsv(locator,nextPageLocator,file)",
    "saveText(var,locator)",
    "saveTextArray(var,locator)",
    "saveTextSubstringAfter(var,locator,delim)",
    "saveTextSubstringBefore(var,locator,delim)",
    "saveTextSubstringBetween(var,locator,start,end)",
    "saveValue(var,locator)",
    "saveValues(var,locator)",
    "screenshot(file,locator)",
    "scrollElement(locator,xOffset,yOffset)",
    "scrollLeft(locator,pixel)",
    "scrollPage(xOffset,yOffset)",
    "scrollRight(locator,pixel)",
    "scrollTo(locator)",
    "select(locator,text)",
    "selectFrame(locator)",
    "selectMulti(locator,array)",
    "selectMultiOptions(locator)",
    "selectText(locator)",
    "selectWindow(winId)",
    "selectWindowAndWait(winId,waitMs)",
    "selectWindowByIndex(index)",
    "selectWindowByIndexAndWait(index,waitMs)",
    "toggleSelections(locator)",
    "type(locator,value)",
    "typeKeys(locator,value)",
    "uncheckAll(locator)",
    "unselectAllText()",
    "updateAttribute(locator,attrName,value)",
    "upload(fieldLocator,file)",
    "verifyContainText(locator,text)",
    "verifyText(locator,text)",
    "wait(waitMs)",
    "waitForElementPresent(locator)",
    "waitForPopUp(winId,waitMs)",
    "waitForTextPresent(text)",
    "waitForTitle(text)"
)
# move existing rows 99-135 down to 101-137 (process bottom-up so we never overwrite
# a source row before it has been read)
for ($i = $webValues.Length - 1; $i -ge 0; $i--) {
    $ws.Cells.Item(101 + $i, 26).Value = $webValues[$i]
}
$ws.Cells.Item(99, 26).Value = "saveSelectedText(var,locator)"
$ws.Cells.Item(100, 26).Value = "saveSelectedValue(var,locator)"

# --- 4. Update the defined names to reflect the new ranges ---
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$39"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$6"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$137"
